$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period label in the remaining data row (2507 -> 2508)
$ws.Range("E16").Value = "2508"

# Update the "Valor Mora" total
$ws.Range("E11").Value = 43333

# Update "Cant. Periodos" count
$ws.Range("F13").Value = 1

# Update the mora value on the single remaining period row
$ws.Range("F16").Value = 43333

# Remove the now-obsolete period rows (old rows 17-23); this shifts
# the signature block (old rows 28-29) up to rows 21-22.
$ws.Rows("17:23").Delete()
